$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Guest State"

# Fix header text/casing for the hotel info table (row 1):
#  - column B is "Floor Number", column C is "Rooms Per Floor"
$ws.Range("B1").Value = "Floor Number"
$ws.Range("C1").Value = "Rooms Per Floor"

# Make row 1 headers use the same bold/centered/wrapped style as the
# room-table headers in row 3 (copy formatting only).
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Give the hotel-name/value row a vertical-centered, wrapped look
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").WrapText = $true

# The floor-count values next to the hotel name pick up the same
# centered/wrapped numeric look used elsewhere in the data rows
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B2:C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Taller header rows to fit wrapped text
$ws.Rows.Item(1).RowHeight = 28.5
$ws.Rows.Item(3).RowHeight = 28.5

# Adjust column widths (closest achievable values to the target
# 19.9296875 / 19.796875 / 19.265625 character widths)
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 19.0
$ws.Columns.Item(3).ColumnWidth = 18.5

# Update selection
$ws.Range("G3").Select() | Out-Null

Write-Host "done"
